# Update the cryptocurrency price/volume snapshot (GitHub Actions refresh).
# Price cells that are plain decimal numbers (e.g. "210.99") are written with
# a leading apostrophe so Excel stores them as text, matching the workbook's
# existing inline-string convention instead of silently converting them to
# floating point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.392.92'
$ws.Range('E2').Value = '  -0.67%  '
$ws.Range('D3').Value = '1.639.05'
$ws.Range('E3').Value = '  -1.59%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '''210.99'
$ws.Range('D6').Value = '''0.532'
$ws.Range('E6').Value = '  +3.66%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').Value = '''23.06'
$ws.Range('E8').Value = '  -2.35%  '
$ws.Range('E9').Value = '  -3.24%  '
$ws.Range('E10').Value = '  -2.24%  '
$ws.Range('D11').Value = '''0.0891'
$ws.Range('E11').Value = '  +1.01%  '
$ws.Range('E12').Value = '  -1.51%  '
$ws.Range('D13').Value = '1.637.37'
$ws.Range('E13').Value = '  -1.40%  '
$ws.Range('E14').Value = '  -2.84%  '
$ws.Range('E15').Value = '  -1.24%  '
$ws.Range('D16').Value = '''64.18'
$ws.Range('E16').Value = '  -2.99%  '
$ws.Range('D17').Value = '27.372.71'
$ws.Range('E17').Value = '  -0.74%  '
$ws.Range('D18').Value = '''229.70'
$ws.Range('E18').Value = '  -5.19%  '
$ws.Range('D19').Value = '0.0₃0719'
$ws.Range('E19').Value = '  -1.50%  '
$ws.Range('D20').Value = '''7.48'
$ws.Range('E20').Value = '  -1.20%  '
$ws.Range('D22').Value = '''4.30'
$ws.Range('E22').Value = '  -4.10%  '
$ws.Range('D23').Value = '''9.32'
$ws.Range('E23').Value = '  +0.39%  '
$ws.Range('D24').Value = '''2.02'
$ws.Range('E24').Value = '  -1.30%  '
$ws.Range('D25').Value = '''147.97'
$ws.Range('E25').Value = '  +1.03%  '
$ws.Range('D26').Value = '''6.94'
$ws.Range('E26').Value = '  -3.47%  '
$ws.Range('E28').Value = '  +0.00%  '
$ws.Range('D29').Value = '''15.51'
$ws.Range('E29').Value = '  -5.32%  '
$ws.Range('E30').Value = '  -5.18%  '
$ws.Range('D31').Value = '''0.0483'
$ws.Range('E31').Value = '  -4.14%  '
$ws.Range('E32').Value = '  -2.38%  '
$ws.Range('E33').Value = '  -0.34%  '
$ws.Range('D34').Value = '1.407.73'
$ws.Range('E34').Value = '  -4.69%  '
$ws.Range('D35').Value = '''1.56'
$ws.Range('E35').Value = '  -0.19%  '
$ws.Range('E36').Value = '  -0.34%  '
$ws.Range('E37').Value = '  -2.20%  '
$ws.Range('D38').Value = '''0.878'
$ws.Range('E38').Value = '  -5.89%  '
$ws.Range('E39').Value = '  -3.16%  '
$ws.Range('E40').Value = '  +1.04%  '
$ws.Range('E41').Value = '  +0.04%  '
$ws.Range('E42').Value = '  -1.40%  '
$ws.Range('E43').Value = '  +1.13%  '
$ws.Range('E44').Value = '  +0.28%  '
$ws.Range('B45').Value = 'TrustWalletToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D45').Value = '''0.789'
$ws.Range('E45').Value = '  -0.06%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').Value = '''64.41'
$ws.Range('E46').Value = '  -7.22%  '
$ws.Range('D47').Value = '1.782.31'
$ws.Range('E48').Value = '  -4.46%  '
$ws.Range('D49').Value = '''87.21'
$ws.Range('E49').Value = '  -2.46%  '
$ws.Range('E50').Value = '  -2.78%  '
$ws.Range('D51').Value = '''0.0987'
$ws.Range('E51').Value = '  -3.76%  '
